$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.041.42'
$ws.Range('E2').Value = '  -0.49%  '

$ws.Range('D3').Value = '1.830.75'
$ws.Range('E3').Value = '  -0.13%  '

$ws.Range('D4').Value = '''0.9987'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '''241.29'
$ws.Range('E5').Value = '  -0.28%  '

$ws.Range('D6').Value = '''0.6233'
$ws.Range('E6').Value = '  -5.19%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '''0.07548'
$ws.Range('E8').Value = '  +2.00%  '

$ws.Range('E9').Value = '  +6.62%  '

$ws.Range('E10').Value = '  -0.67%  '

$ws.Range('D11').Value = '''22.82'
$ws.Range('E11').Value = '  -0.10%  '

$ws.Range('D12').Value = '''0.07636'
$ws.Range('E12').Value = '  -1.67%  '

$ws.Range('D13').Value = '1.828.18'
$ws.Range('E13').Value = '  -0.83%  '

$ws.Range('D14').Value = '''4.959'
$ws.Range('E14').Value = '  -0.73%  '

$ws.Range('D15').Value = '''0.6653'
$ws.Range('E15').Value = '  -0.12%  '

$ws.Range('D16').Value = '''82.34'

$ws.Range('D17').Value = '''0.000009100'
$ws.Range('E17').Value = '  +8.29%  '

$ws.Range('D18').Value = '''6.001'
$ws.Range('E18').Value = '  -1.66%  '

$ws.Range('D19').Value = '29.035.55'
$ws.Range('E19').Value = '  -0.44%  '

$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.084.89'
$ws.Range('E20').Value = '  +1.05%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''225.24'
$ws.Range('E21').Value = '  -0.82%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '''12.35'
$ws.Range('E22').Value = '  -0.80%  '

$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '''7.186'
$ws.Range('E24').Value = '  +0.84%  '

$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = '''1.0000'
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''159.75'
$ws.Range('E26').Value = '  +0.41%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''8.392'
$ws.Range('E27').Value = '  -2.49%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '''0.1358'
$ws.Range('E28').Value = '  -2.10%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''17.84'
$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.495'
$ws.Range('E30').Value = '  -1.46%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''4.055'
$ws.Range('E31').Value = '  -1.38%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''1.205'
$ws.Range('E32').Value = '  +0.89%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''4.027'
$ws.Range('E33').Value = '  -0.41%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.05206'
$ws.Range('E34').Value = '  -1.25%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '''1.837'
$ws.Range('E35').Value = '  -1.20%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''1.155'
$ws.Range('E36').Value = '  +1.37%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.7315'
$ws.Range('E37').Value = '  -1.16%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '''2.602'
$ws.Range('E38').Value = '  -1.99%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.284.31'
$ws.Range('E39').Value = '  -1.28%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.763'
$ws.Range('E40').Value = '  +1.15%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.01779'
$ws.Range('E41').Value = '  -0.64%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''6.369'
$ws.Range('E42').Value = '  +7.45%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.8914'
$ws.Range('E43').Value = '  -3.81%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''101.42'
$ws.Range('E45').Value = '  -0.92%  '

$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.980.90'
$ws.Range('E46').Value = '  +0.87%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '''0.5114'
$ws.Range('E47').Value = '  -0.55%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''63.51'
$ws.Range('E48').Value = '  +1.02%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '''0.00000000119'
$ws.Range('E49').Value = '  -0.51%  '

$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '''0.3966'
$ws.Range('E50').Value = '  -1.02%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''8.864'
$ws.Range('E51').Value = '  +1.05%  '
